$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 10372.65132737054
$ws2025.Range("E2").Value = 289260.5393052954
$ws2025.Range("G2").Value = 80959.25712661834
$ws2025.Range("I2").Value = 161710.6685703679
$ws2025.Range("L2").Value = 484922.2142001599
$ws2025.Range("M2").Value = 105953.7713982
$ws2025.Range("N2").Value = 70003.73489578845
$ws2025.Range("O2").Value = 69744.89343456978

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 0
$ws2030.Range("B2").Value = 31203.23858116339
$ws2030.Range("E2").Value = 170658.5511254234
$ws2030.Range("I2").Value = 209080.6134235085
$ws2030.Range("L2").Value = 63518.11613148725
$ws2030.Range("M2").Value = 68536.72857011756
$ws2030.Range("N2").Value = 19285.19160463996
$ws2030.Range("O2").Value = 27033.1386905727

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 27543.1755456332
$ws2035.Range("B2").Value = 22113.21643273498
$ws2035.Range("E2").Value = 114655.4402706629
$ws2035.Range("I2").Value = 153866.0861464091
$ws2035.Range("M2").Value = 44638.22942194272
$ws2035.Range("N2").Value = 39676.88529639924
$ws2035.Range("O2").Value = 31311.04369977792

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("N2").Value = 1142.580190039942
$ws2040.Range("O2").Value = 0

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 29588.33508286276
$ws2045.Range("N2").Value = 4347.543515635315
$ws2045.Range("O2").Value = 20429.76977394434
